$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (even_MAG-GUT12458.fa and even_MAG-GUT20500.fa),
# shifting rows 4 and 5 up to become the new rows 2 and 3.
$ws.Rows("2:3").Delete()
